$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "945÷8=118, 1" "913÷3=304, 1"
Replace-Text "844÷4=211, 0" "376÷2=188, 0"
Replace-Text "705÷9=78, 3" "557÷6=92, 5"
Replace-Text "954÷5=190, 4" "951÷5=190, 1"
Replace-Text "278÷8=34, 6" "941÷9=104, 5"
Replace-Text "252÷8=31, 4" "681÷8=85, 1"
Replace-Text "366÷4=91, 2" "231÷6=38, 3"
Replace-Text "322÷5=64, 2" "159÷8=19, 7"
Replace-Text "934÷8=116, 6" "211÷7=30, 1"
Replace-Text "800÷9=88, 8" "591÷7=84, 3"
Replace-Text "556÷2=278, 0" "829÷2=414, 1"
Replace-Text "164÷4=41, 0" "345÷6=57, 3"
Replace-Text "594÷9=66, 0" "930÷2=465, 0"
Replace-Text "748÷4=187, 0" "835÷7=119, 2"
Replace-Text "317÷5=63, 2" "849÷5=169, 4"
Replace-Text "666÷5=133, 1" "956÷2=478, 0"
Replace-Text "413÷5=82, 3" "170÷3=56, 2"
Replace-Text "961÷7=137, 2" "124÷3=41, 1"
Replace-Text "994÷2=497, 0" "209÷8=26, 1"
Replace-Text "253÷3=84, 1" "716÷2=358, 0"
Replace-Text "463÷9=51, 4" "814÷7=116, 2"
Replace-Text "470÷2=235, 0" "959÷2=479, 1"
Replace-Text "426÷3=142, 0" "365÷7=52, 1"
Replace-Text "373÷3=124, 1" "805÷9=89, 4"
Replace-Text "847÷6=141, 1" "510÷3=170, 0"
